# 22: temp commit: sign agreement
#
# Sheet "Confirmation Events": add a new "Sign Agreement" row (row 4) with a
# due date, styled with a small Menlo font/color, and widen column A so the
# new label fits.
#
# Sheet "Candidates with events": add a third confirmed-events pair of
# columns (candidate_events.2.completed_date / candidate_events.2.admin_confirmed)
# with sample data for the first two candidate rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Confirmation Events"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Widen column A to fit the new "Sign Agreement" label.
$ws1.Columns.Item(1).ColumnWidth = 24.43

$nameCell = $ws1.Range("A4")
$nameCell.Value = "Sign Agreement"
$nameCell.Font.Name = "Menlo"
$nameCell.Font.Size = 10
$nameCell.Font.Color = 197380

$dateCell = $ws1.Range("B4")
$dateCell.Value = 42564
$dateCell.NumberFormat = "m/d/yyyy"

# ---------------------------------------------------------------------------
# Sheet 2: "Candidates with events"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Insert 3 new columns (R:T) for the 3rd candidate event (completed_date,
# admin_confirmed, and a trailing blank spacer column matching the sheet's
# existing layout pattern).
$ws2.Range("R1:T1").EntireColumn.Insert()
$ws2.Range("R1:T1").EntireColumn.ColumnWidth = 23.43

$ws2.Range("R1").Value = "candidate_events.2.completed_date"
$ws2.Range("S1").Value = "candidate_events.2.admin_confirmed"

$ws2.Range("R2").Value = 42441
$ws2.Range("R2").NumberFormat = "m/d/yyyy"
$ws2.Range("S2").Value = $true

Write-Output "edit applied"
